$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# ARKCORR-25: Fixing days in queue
#
# The "Set Field Value" column that duplicated the queue-enter-date action
# (old column E) is removed - it was an exact duplicate of column D and the
# real per-queue "owning group" actions lived one column over (old column F,
# which becomes the new column E after the delete).
# ---------------------------------------------------------------------------

# Remove the old (duplicate) column E -- shifts old F => new E.
$ws.Columns.Item(5).Delete()

# Row 21's "When Expression is True" cell (D21) picks up the thin right
# border that the rest of the table row uses now that the table is one
# column narrower (style index 4 instead of 10).
$ws.Range("D21").Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------------------
# New rules: track the previous queue so "queue enter date" rules only fire
# when the case actually transitions INTO that queue (not on every save),
# plus a "Release" queue-enter-date rule and the "Set Previous Queue" rule
# that records the queue before it changes.
# ---------------------------------------------------------------------------

# Use existing, fully-styled rows as a formatting template for the new rows
# so fills/borders/fonts match the rest of the rule table.
$ws.Range("A27:E27").Copy()
$ws.Range("A35:E38").PasteSpecial(-4122)
$ws.Range("A27:E27").Copy()
$ws.Range("A39:E40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B35").Value = "Set Queue Enter Date Intake Queue"
$ws.Range("C35").Value = 'queue?.name == "Intake" && (previousQueue == null || previousQueue?.name != "Intake")'
$ws.Range("D35").Value = "setQueueEnterDate, java.time.LocalDate.now()"

$ws.Range("B36").Value = "Set Queue Enter Date Fulfill Queue"
$ws.Range("C36").Value = 'queue?.name == "Fulfill" && (previousQueue == null || previousQueue?.name != "Fulfill")'
$ws.Range("D36").Value = "setQueueEnterDate, java.time.LocalDate.now()"

$ws.Range("B37").Value = "Set Queue Enter Date Supervisor Approval Queue"
$ws.Range("C37").Value = 'queue?.name == "Supervisor Approval" && (previousQueue == null || previousQueue?.name != "Supervisor Approval")'
$ws.Range("D37").Value = "setQueueEnterDate, java.time.LocalDate.now()"

$ws.Range("B38").Value = "Set Queue Enter Date Executive Approval Queue"
$ws.Range("C38").Value = 'queue?.name == "Executive Approval" && (previousQueue == null || previousQueue?.name != "Executive Approval")'
$ws.Range("D38").Value = "setQueueEnterDate, java.time.LocalDate.now()"

$ws.Range("B39").Value = "Set Queue Enter Date Intake Queue"
$ws.Range("C39").Value = 'queue?.name == "Release" && (previousQueue == null || previousQueue?.name != "Release")'
$ws.Range("D39").Value = "setQueueEnterDate, java.time.LocalDate.now()"

$ws.Range("B40").Value = "Set Previous Queue"
$ws.Range("C40").Value = "queue == null || queue != null"
$ws.Range("D40").Value = "setPreviousQueue, `$caseFile.getQueue()"

# ---------------------------------------------------------------------------
# Column widths (best-fit to the new, longer condition/action text).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.42578125
$ws.Columns.Item(2).ColumnWidth = 45.85546875
$ws.Columns.Item(3).ColumnWidth = 109.5703125
$ws.Columns.Item(4).ColumnWidth = 100.42578125
$ws.Columns.Item(5).ColumnWidth = 44

# Move the active selection/view to where the new rule rows were added.
$ws.Activate()
$ws.Range("D40").Select()
